$d = $word.ActiveDocument

# --- Step 1: Remove the "License Information" heading paragraph (Heading2) ---
$licenseHeadingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("License Information")) {
        $licenseHeadingPara = $p
        break
    }
}
if ($licenseHeadingPara -ne $null) {
    $licenseHeadingPara.Range.Delete()
}

# --- Step 2: Locate the resource-description paragraph (the one that currently
#     starts with "Translation Questions (unfoldingWord) is based on") and the
#     paragraph right after it ("This PDF version..."). We merge the two into a
#     single paragraph and replace the text with the new resource description. ---
$descPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Translation Questions (unfoldingWord) is based on")) {
        $descPara = $p
        break
    }
}

$followingPara = $descPara.Next()

# Clear the text of each paragraph individually (same-paragraph ranges only,
# since deleting a range that spans into the middle of a following paragraph's
# text is unreliable). Then merge the two (now-empty) paragraphs into one by
# deleting the first paragraph's mark.
$delText1 = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$delText1.Delete()

$followingPara2 = $d.Paragraphs.Item($descPara.Range.Start + 1)
$delText2 = $d.Range($followingPara2.Range.Start, $followingPara2.Range.End - 1)
$delText2.Delete()

$mergedStart = $descPara.Range.Start
$emptyPara = $d.Paragraphs.Item($mergedStart + 1)
$markRange = $d.Range($emptyPara.Range.Start, $emptyPara.Range.End)
$markRange.Delete()

# --- Step 3: Insert the new paragraph content ---
$cur = $mergedStart

$boldText = "unfoldingWord® Translation Questions"
$seg2 = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. "
$seg3 = "unfoldingWord® Translation Questions"
$seg4 = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from "
$seg5 = "unfoldingWord® Translation Questions"
$seg6 = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

$ins = $d.Range($cur, $cur)
$ins.InsertAfter($boldText)
$newEnd = $cur + $boldText.Length
$segRng = $d.Range($cur, $newEnd)
$segRng.Font.Bold = $true
$cur = $newEnd

$ins = $d.Range($cur, $cur)
$ins.InsertAfter($seg2)
$newEnd = $cur + $seg2.Length
$segRng = $d.Range($cur, $newEnd)
$segRng.Font.Bold = $false
$cur = $newEnd

$ins = $d.Range($cur, $cur)
$ins.InsertAfter($seg3)
$newEnd = $cur + $seg3.Length
$segRng = $d.Range($cur, $newEnd)
$segRng.Font.Bold = $false
$cur = $newEnd

$ins = $d.Range($cur, $cur)
$ins.InsertAfter($seg4)
$newEnd = $cur + $seg4.Length
$segRng = $d.Range($cur, $newEnd)
$segRng.Font.Bold = $false
$cur = $newEnd

$ins = $d.Range($cur, $cur)
$ins.InsertAfter($seg5)
$newEnd = $cur + $seg5.Length
$segRng = $d.Range($cur, $newEnd)
$segRng.Font.Bold = $false
$cur = $newEnd

$ins = $d.Range($cur, $cur)
$ins.InsertAfter($seg6)
$newEnd = $cur + $seg6.Length
$segRng = $d.Range($cur, $newEnd)
$segRng.Font.Bold = $false
$cur = $newEnd

Write-Host "Edit complete"
